$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update status of row 2 (Health Status of Player) from "Started" to "Done"
$ws.Range("E2").Value = "Done"

# Add status "Started" to row 13 (Main Menu)
$ws.Range("E13").Value = "Started"

# Add status "Done" to row 15 (Should not fall off the map)
$ws.Range("E15").Value = "Done"

# Update the active selection on the sheet to F13
$ws.Range("F13").Select()
